{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n// (percentages, dollar amounts, large numbers) across the resume body.\n//\n// For each target paragraph we locate the metric substring(s) via\n// paragraph.search() and set font.bold = true plus font.color = \"#2C3E50\"\n// on the matched range(s). Word automatically splits the run(s) so that\n// only the matched text carries the new formatting, leaving the rest of\n// the paragraph's text in plain runs \u2014 mirroring the OOXML diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Map of unique paragraph text (exact match) -> list of metric substrings\n// that must be bolded + colored, in the order they appear.\nconst targets = [\n  {\n    text: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"]\n  },\n  {\n    text: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n    metrics: [\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\"]\n  },\n  {\n    text: \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"]\n  },\n  {\n    text: \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"]\n  },\n  {\n    text: \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"]\n  },\n  {\n    text: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"]\n  }\n];\n\n// Track how many times we've matched each distinct target text, so that\n// duplicate paragraph texts (none in this document, but kept robust)\n// are each only processed once via their own paragraph object.\nconst usedParagraphIndices = new Set();\n\nfor (const target of targets) {\n  let matchedIndex = -1;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (usedParagraphIndices.has(i)) continue;\n    if (paragraphs.items[i].text === target.text) {\n      matchedIndex = i;\n      break;\n    }\n  }\n  if (matchedIndex === -1) {\n    throw new Error(\"Could not find target paragraph: \" + target.text);\n  }\n  usedParagraphIndices.add(matchedIndex);\n  const para = paragraphs.items[matchedIndex];\n\n  for (const metric of target.metrics) {\n    const results = para.search(metric, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    if (results.items.length === 0) {\n      throw new Error(\"Metric not found in paragraph: \" + metric);\n    }\n    const range = results.items[0];\n    range.font.bold = true;\n    range.font.color = HIGHLIGHT_COLOR;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n# (percentages, dollar amounts, large numbers) across the resume body.\n#\n# For each target paragraph, we walk its list of metric substrings in\n# order, using Find.Execute scoped to a shrinking sub-range (from the end\n# of the previous match through the paragraph end) so repeated/overlapping\n# substrings are matched left-to-right without re-matching earlier text.\n# Each match's Font.Bold / Font.Color is then set, which causes Word to\n# split the run so only the metric text carries the new formatting.\n\n$d = $word.ActiveDocument\n\n# Target color #2C3E50 expressed as the BGR-packed decimal Word's Font.Color\n# (wdColor) expects: val = (B << 16) | (G << 8) | R.\n#   R=0x2C=44, G=0x3E=62, B=0x50=80  ->  (80*65536)+(62*256)+44 = 5258796\n$highlightColor = 5258796\n\nfunction Format-Metrics {\n    param([object]$Paragraph, [object]$Metrics)\n    $pStart = $Paragraph.Range.Start\n    $pEnd = $Paragraph.Range.End\n    $searchStart = $pStart\n    foreach ($metric in $Metrics) {\n        $r = $d.Range($searchStart, $pEnd)\n        $found = $r.Find.Execute($metric)\n        if (-not $found) {\n            throw \"Metric '$metric' not found in paragraph starting at $pStart\"\n        }\n        $r.Font.Bold = 1\n        $r.Font.Color = $highlightColor\n        $searchStart = $r.End\n    }\n}\n\nforeach ($p in $d.Paragraphs) {\n    # Range.Text includes the trailing paragraph mark (CR, char 13) \u2014 and\n    # occasionally a cell/section mark (char 7) \u2014 so strip those before\n    # comparing against literal target strings.\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($t -eq \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\") {\n        Format-Metrics $p @(\"23%\", \"64%\")\n    }\n    elseif ($t -eq \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\") {\n        Format-Metrics $p @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\")\n    }\n    elseif ($t -eq \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\") {\n        Format-Metrics $p @(\"1,200\")\n    }\n    elseif ($t -eq \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\") {\n        Format-Metrics $p @(\"`$400M\", \"`$1B\")\n    }\n    elseif ($t -eq \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\") {\n        Format-Metrics $p @(\"73.5%\", \"`$4.7M\")\n    }\n    elseif ($t -eq \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\") {\n        Format-Metrics $p @(\"87%\", \"71%\")\n    }\n}\n"}
